$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Numéro de page" value in row 123 (this script run found nothing
# relevant, so the page-number column is blank, just like the other
# "Rien ne nous concerne aujourd'hui !" rows)
$ws.Range("C123").Value = ""
$ws.Range("C123").Style = "Normal"

# Append a new row 124 with today's script result
# Force A124 to be stored as text so "2025-06-04" is not auto-converted
# into a date serial number
$ws.Range("A124").NumberFormat = "@"
$ws.Range("A124").Value = "2025-06-04"
$ws.Range("A124").Style = "Normal"

$ws.Range("B124").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C124").Value = "NA"
$ws.Range("D124").Value = 1
